# feat: update obby new level
# Appends 5 new checkpoint rows (30-34) to the "Obbycheck" sheet, mirroring
# the existing id/checkpointloc/splashpointloc/splashscale layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{ Id = 30; CheckpointLoc = "393397.53|8672.15|31400.97"; SplashPointLoc = "393400.94|8666.14|31052.23" },
    @{ Id = 31; CheckpointLoc = "396374.53|8289.15|31400.97"; SplashPointLoc = "396378.94|8299.14|31052.23" },
    @{ Id = 32; CheckpointLoc = "394144.53|7755.15|31599.97"; SplashPointLoc = "394154.94|7762.14|31500.23" },
    @{ Id = 33; CheckpointLoc = "396558.53|8610.15|32502.97"; SplashPointLoc = "396563.94|8594.14|32177.23" },
    @{ Id = 34; CheckpointLoc = "399253.53|8610.15|33974.97"; SplashPointLoc = "399263.94|8611.14|33564.23" }
)

$startRow = 35
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $row = $startRow + $i
    $data = $newRows[$i]
    $ws.Range("A$row").Value = $data.Id
    $ws.Range("B$row").Value = $data.CheckpointLoc
    $ws.Range("C$row").Value = $data.SplashPointLoc
    $ws.Range("D$row").Value = "2.0|2.0|2.0"
}

# Move the active selection to the last newly added cell, as in the authored edit.
$ws.Range("D39").Select()
